# Applies the scheduled-runner price/profit refresh to the Halicarnassus_Profits
# workbook: updated currentAveragePrice* / LevePrice* / LeveProfit* figures across
# the eight job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1624.5
$ws.Range("I19").Value = 1833
$ws.Range("J19").Value = 999
$ws.Range("K19").Value = 1833
$ws.Range("L19").Value = 999
$ws.Range("M19").Value = -1658
$ws.Range("H32").Value = 697.75
$ws.Range("I32").Value = 501
$ws.Range("J32").Value = 763.3333
$ws.Range("K32").Value = 501
$ws.Range("L32").Value = 763.3333
$ws.Range("M32").Value = -175
$ws.Range("N32").Value = -1415.3333
$ws.Range("H70").Value = 4572.364
$ws.Range("I70").Value = 3190.4546
$ws.Range("J70").Value = 5954.273
$ws.Range("K70").Value = 9571.363799999999
$ws.Range("L70").Value = 17862.819
$ws.Range("M70").Value = -9301.363799999999
$ws.Range("N70").Value = -18402.819
$ws.Range("H73").Value = 4572.364
$ws.Range("I73").Value = 3190.4546
$ws.Range("J73").Value = 5954.273
$ws.Range("K73").Value = 9571.363799999999
$ws.Range("L73").Value = 17862.819
$ws.Range("M73").Value = -8635.363799999999
$ws.Range("N73").Value = -19734.819
$ws.Range("H86").Value = 2752
$ws.Range("I86").Value = 2500
$ws.Range("J86").Value = 3004
$ws.Range("K86").Value = 2500
$ws.Range("L86").Value = 3004
$ws.Range("M86").Value = -1377
$ws.Range("N86").Value = -5250
$ws.Range("H89").Value = 2752
$ws.Range("I89").Value = 2500
$ws.Range("J89").Value = 3004
$ws.Range("K89").Value = 12500
$ws.Range("L89").Value = 15020
$ws.Range("M89").Value = -6884
$ws.Range("N89").Value = -26252
$ws.Range("H132").Value = 10641.407
$ws.Range("I132").Value = 10046.318
$ws.Range("J132").Value = 13259.8
$ws.Range("K132").Value = 30138.954
$ws.Range("L132").Value = 39779.39999999999
$ws.Range("M132").Value = -27608.954
$ws.Range("N132").Value = -44839.39999999999
$ws.Range("H137").Value = 2508.5454
$ws.Range("I137").Value = 1099
$ws.Range("J137").Value = 4200
$ws.Range("K137").Value = 3297
$ws.Range("L137").Value = 12600
$ws.Range("M137").Value = -747

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 14732.4
$ws.Range("I28").Value = 14732.4
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 14732.4
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -14540.4
$ws.Range("H61").Value = 2018.2354
$ws.Range("I61").Value = 1434.6666
$ws.Range("J61").Value = 3418.8
$ws.Range("K61").Value = 1434.6666
$ws.Range("L61").Value = 3418.8
$ws.Range("M61").Value = -1222.6666
$ws.Range("N61").Value = -3842.8
$ws.Range("H99").Value = 14732.4
$ws.Range("I99").Value = 14732.4
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 14732.4
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -11737.4
$ws.Range("H102").Value = 9620649
$ws.Range("I102").Value = 17858730
$ws.Range("J102").Value = 9553.166999999999
$ws.Range("K102").Value = 17858730
$ws.Range("L102").Value = 9553.166999999999
$ws.Range("M102").Value = -17857108
$ws.Range("N102").Value = -12797.167
$ws.Range("H118").Value = 80000
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 80000
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 80000
$ws.Range("N118").Value = -83314
$ws.Range("H136").Value = 2018.2354
$ws.Range("I136").Value = 1434.6666
$ws.Range("J136").Value = 3418.8
$ws.Range("K136").Value = 4303.9998
$ws.Range("L136").Value = 10256.4
$ws.Range("M136").Value = -1753.9998
$ws.Range("N136").Value = -15356.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 21427.416
$ws.Range("I82").Value = 10695.8
$ws.Range("J82").Value = 75085.5
$ws.Range("K82").Value = 10695.8
$ws.Range("L82").Value = 75085.5
$ws.Range("M82").Value = -10312.8
$ws.Range("N82").Value = -75851.5
$ws.Range("H85").Value = 21427.416
$ws.Range("I85").Value = 10695.8
$ws.Range("J85").Value = 75085.5
$ws.Range("K85").Value = 10695.8
$ws.Range("L85").Value = 75085.5
$ws.Range("M85").Value = -9369.799999999999
$ws.Range("N85").Value = -77737.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 319.66666
$ws.Range("I12").Value = 300
$ws.Range("J12").Value = 329.5
$ws.Range("K12").Value = 300
$ws.Range("L12").Value = 329.5
$ws.Range("M12").Value = -130
$ws.Range("N12").Value = -669.5
$ws.Range("H22").Value = 2267.111
$ws.Range("I22").Value = 1254.25
$ws.Range("J22").Value = 3077.4
$ws.Range("K22").Value = 1254.25
$ws.Range("L22").Value = 3077.4
$ws.Range("M22").Value = -904.25
$ws.Range("N22").Value = -3777.4
$ws.Range("H31").Value = 3596.2273
$ws.Range("I31").Value = 1624.3334
$ws.Range("J31").Value = 4961.385
$ws.Range("K31").Value = 1624.3334
$ws.Range("L31").Value = 4961.385
$ws.Range("M31").Value = -1329.3334
$ws.Range("N31").Value = -5551.385
$ws.Range("H34").Value = 3596.2273
$ws.Range("I34").Value = 1624.3334
$ws.Range("J34").Value = 4961.385
$ws.Range("K34").Value = 1624.3334
$ws.Range("L34").Value = 4961.385
$ws.Range("M34").Value = -1422.3334
$ws.Range("N34").Value = -5365.385
$ws.Range("H86").Value = 3375
$ws.Range("I86").Value = 3375
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3375
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2252
$ws.Range("H89").Value = 3375
$ws.Range("I89").Value = 3375
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 16875
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -11259
$ws.Range("H132").Value = 1892.2941
$ws.Range("I132").Value = 1892.2941
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5676.8823
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3146.8823
$ws.Range("H134").Value = 2226.2144
$ws.Range("I134").Value = 2226.2144
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6678.6432
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -4143.6432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 161.13333
$ws.Range("I12").Value = 42.25
$ws.Range("J12").Value = 204.36363
$ws.Range("K12").Value = 126.75
$ws.Range("L12").Value = 613.0908899999999
$ws.Range("M12").Value = 46.25
$ws.Range("N12").Value = -959.0908899999999
$ws.Range("H131").Value = 2196.3572
$ws.Range("I131").Value = 1287.5
$ws.Range("J131").Value = 2559.9
$ws.Range("K131").Value = 3862.5
$ws.Range("L131").Value = 7679.700000000001
$ws.Range("M131").Value = 1177.5
$ws.Range("N131").Value = -17759.7
$ws.Range("H140").Value = 2429.1333
$ws.Range("I140").Value = 1649.8182
$ws.Range("J140").Value = 4572.25
$ws.Range("K140").Value = 4949.4546
$ws.Range("L140").Value = 13716.75
$ws.Range("M140").Value = 230.5454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 61613
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 61613
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 61613
$ws.Range("N57").Value = -63253
$ws.Range("H122").Value = 1875.5294
$ws.Range("I122").Value = 1197
$ws.Range("J122").Value = 2350.5
$ws.Range("K122").Value = 3591
$ws.Range("L122").Value = 7051.5
$ws.Range("M122").Value = -1141
$ws.Range("N122").Value = -11951.5
$ws.Range("H126").Value = 420
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 420
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 1260
$ws.Range("N126").Value = -6200
$ws.Range("H132").Value = 11919
$ws.Range("I132").Value = 14448.75
$ws.Range("J132").Value = 1800
$ws.Range("K132").Value = 43346.25
$ws.Range("L132").Value = 5400
$ws.Range("M132").Value = -40816.25
$ws.Range("N132").Value = -10460

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 500
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 500
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -330
$ws.Range("H22").Value = 1007.53845
$ws.Range("I22").Value = 850
$ws.Range("J22").Value = 1077.5555
$ws.Range("K22").Value = 850
$ws.Range("L22").Value = 1077.5555
$ws.Range("M22").Value = -555
$ws.Range("N22").Value = -1667.5555
$ws.Range("H27").Value = 1007.53845
$ws.Range("I27").Value = 850
$ws.Range("J27").Value = 1077.5555
$ws.Range("K27").Value = 850
$ws.Range("L27").Value = 1077.5555
$ws.Range("M27").Value = -743
$ws.Range("N27").Value = -1291.5555
$ws.Range("H55").Value = 1367
$ws.Range("I55").Value = 1635.8334
$ws.Range("J55").Value = 1165.375
$ws.Range("K55").Value = 1635.8334
$ws.Range("L55").Value = 1165.375
$ws.Range("M55").Value = -1462.8334
$ws.Range("N55").Value = -1511.375
$ws.Range("H82").Value = 2939.0625
$ws.Range("I82").Value = 634.2857
$ws.Range("J82").Value = 4731.6665
$ws.Range("K82").Value = 634.2857
$ws.Range("L82").Value = 4731.6665
$ws.Range("M82").Value = -273.2857
$ws.Range("N82").Value = -5453.6665
$ws.Range("H85").Value = 2939.0625
$ws.Range("I85").Value = 634.2857
$ws.Range("J85").Value = 4731.6665
$ws.Range("K85").Value = 634.2857
$ws.Range("L85").Value = 4731.6665
$ws.Range("M85").Value = 613.7143
$ws.Range("N85").Value = -7227.6665
$ws.Range("H93").Value = 4049.2
$ws.Range("I93").Value = 4049.2
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 4049.2
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -2801.2
$ws.Range("N93").ClearContents()
$ws.Range("H122").Value = 3262.889
$ws.Range("I122").Value = 3370.125
$ws.Range("J122").Value = 2405
$ws.Range("K122").Value = 10110.375
$ws.Range("L122").Value = 7215
$ws.Range("M122").Value = -7660.375
$ws.Range("N122").Value = -12115
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 683.3333
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 683.3333
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 683.3333
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -1141.3333
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H105").Value = 25447
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 25447
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 25447
$ws.Range("N105").Value = -32435
$ws.Range("H118").Value = 100000
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 100000
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 100000
$ws.Range("N118").Value = -103314

